# Apply the cryptocurrency price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "57.119.53"
$ws.Cells.Item(2, 5).Value = "  -1.73%  "

$ws.Cells.Item(3, 4).Value = "2.990.62"
$ws.Cells.Item(3, 5).Value = "  -2.71%  "

$ws.Cells.Item(4, 5).Value = "  +0.14%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "499.96"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -5.38%  "

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "136.16"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -5.21%  "

$ws.Cells.Item(7, 5).Value = "  +0.18%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.428"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -4.58%  "

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.24"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -5.75%  "

$ws.Cells.Item(10, 5).Value = "  -5.83%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.354"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -4.79%  "

$ws.Cells.Item(12, 4).Value = "3.507.32"
$ws.Cells.Item(12, 5).Value = "  -2.42%  "

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.126"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.65%  "

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "25.89"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -5.84%  "

$ws.Cells.Item(15, 5).Value = "  -7.63%  "

$ws.Cells.Item(16, 4).Value = "57.215.01"
$ws.Cells.Item(16, 5).Value = "  -1.53%  "

$ws.Cells.Item(17, 2).Value = "WrappedEther"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(17, 4).Value = "3.014.22"
$ws.Cells.Item(17, 5).Value = "  -1.74%  "

$ws.Cells.Item(18, 2).Value = "Polkadot"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.04"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -3.38%  "

$ws.Cells.Item(19, 5).Value = "  -4.54%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.81"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -4.72%  "

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "318.82"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -6.78%  "

$ws.Cells.Item(22, 5).Value = "  -0.03%  "

$ws.Cells.Item(23, 5).Value = "  +0.46%  "

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.490"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -2.83%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "63.11"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -3.08%  "

$ws.Cells.Item(26, 5).Value = "  -0.13%  "

$ws.Cells.Item(27, 5).Value = "  -5.50%  "

$ws.Cells.Item(28, 4).Value = "0.0₃0879"
$ws.Cells.Item(28, 5).Value = "  -10.86%  "

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.59"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -5.96%  "

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.03"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -5.48%  "

$ws.Cells.Item(31, 5).Value = "  -4.94%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.15"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -7.65%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "20.07"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -5.19%  "

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "154.98"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -1.25%  "

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.51"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -6.03%  "

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.74"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -5.15%  "

$ws.Cells.Item(37, 5).Value = "  -8.45%  "

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "24.12"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -9.42%  "

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0658"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -7.14%  "

$ws.Cells.Item(40, 4).Value = "3.027.69"
$ws.Cells.Item(40, 5).Value = "  -2.63%  "

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "37.84"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.13%  "

$ws.Cells.Item(42, 5).Value = "  -0.01%  "

$ws.Cells.Item(43, 2).Value = "Mantle"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.645"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -3.82%  "

$ws.Cells.Item(44, 2).Value = "Filecoin"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.70"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -5.76%  "

$ws.Cells.Item(45, 4).Value = "2.169.61"
$ws.Cells.Item(45, 5).Value = "  -7.13%  "

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.37"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -8.16%  "

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.91"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.29%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.928"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -10.44%  "

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0232"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -5.84%  "

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.15"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -5.96%  "

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.75"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -13.62%  "

